$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Genre -> nonfiction, Volumes -> 1 (keep as text like the rest of the sheet)
$ws.Range("B2").Value = "nonfiction"
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "1"

# Row 3: becomes what used to be row 4 (dragon ball / manga / 1,2,3,4,5,6,7)
$ws.Range("A3").Value = "dragon ball"
$ws.Range("B3").Value = "manga"
$ws.Range("C3").Value = "1,2,3,4,5,6,7"

# Row 4: becomes what used to be row 5 (naruto / manga / 2)
$ws.Range("A4").Value = "naruto"
$ws.Range("B4").Value = "manga"
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "2"

# Row 5: becomes what used to be row 3 (d / manga / 2)
$ws.Range("A5").Value = "d"
$ws.Range("B5").Value = "manga"
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "2"
